$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.511.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.487.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.57%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.66%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.57"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.79%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.485.86"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.50%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.84%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.25%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.46%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.091.07"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.50%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.14"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.97%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.565.32"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.92%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.488.09"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.78%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.98"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +4.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "390.94"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.90"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.63%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +3.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.18%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.66%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.13%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.33"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +8.92%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.66%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.30"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.80%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.10%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.55"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.64%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.40"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.49%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +8.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.86"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.91"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.78"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.48"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.87%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +6.68%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.11"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.768.54"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0312"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.84%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.34%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "345.22"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.67%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.46%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +12.34%  "
